$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in Member (A) / Task Description (B) columns for rows 2-10 ---
# Written row by row, column A before column B, to mirror the natural
# data-entry order used when this sheet was authored.
$ws.Range("A2").Value = "All"
$ws.Range("B2").Value = "Add more records to the database"

$ws.Range("A3").Value = "Arpit "
$ws.Range("B3").Value = "Make login page"

$ws.Range("A4").Value = "Arpit "
$ws.Range("B4").Value = "Make registration page"

$ws.Range("A5").Value = "Ishika"
$ws.Range("B5").Value = "Make forgot password page"

$ws.Range("A6").Value = "Ishika"
$ws.Range("B6").Value = "Design report on billing"

$ws.Range("A7").Value = "Tri"
$ws.Range("B7").Value = "Make upcoming trips "

$ws.Range("A8").Value = "Maeve"
$ws.Range("B8").Value = "Delete reservations"

$ws.Range("A9").Value = "Maeve"
$ws.Range("B9").Value = "Edit reservations"

$ws.Range("A10").Value = "Tri"
$ws.Range("B10").Value = "Test the pages"

# --- Match column A/B formatting to the existing data-row style for all
#     rows 2-10 (rows 8-10 are brand new and don't yet carry the style) ---
$ws.Range("B2").Copy()
$ws.Range("A2:B10").PasteSpecial(-4122)

# --- Widen column B to fit the new member names / task text ---
$ws.Columns.Item(2).ColumnWidth = 33.6
